$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26").Value = 2400.2938700802683
$ws.Range("E26").Value = 704.228828727142
$ws.Range("F26").Value = 146.2989190295223
$ws.Range("G26").Value = 1196.675294936432
$ws.Range("K26").Value = 0.03325044127125001
$ws.Range("L26").Value = 31358.0309218639
$ws.Range("M26").Value = 1042.6683655494774

$ws.Range("C27").Value = 1316.8943799147532
$ws.Range("E27").Value = 335.3856067465372
$ws.Range("F27").Value = 68.88202738791833
$ws.Range("G27").Value = 726.6255744620322
$ws.Range("K27").Value = 0.027450779798749995
$ws.Range("L27").Value = 21720.636897403987
$ws.Range("M27").Value = 596.2484205592411

$ws.Range("C28").Value = 1370.7322634204334
$ws.Range("E28").Value = 333.450805202819
$ws.Range("F28").Value = 54.2451454382689
$ws.Range("G28").Value = 906.3147735712006
$ws.Range("K28").Value = 0.03250710120375
$ws.Range("L28").Value = 24271.490366281625
$ws.Range("M28").Value = 788.9957937025599

$ws.Range("C29").Value = 2087.5735007623325
$ws.Range("E29").Value = 670.763347186695
$ws.Range("F29").Value = 131.89723230104508
$ws.Range("G29").Value = 989.8496301537058
$ws.Range("K29").Value = 0.035086588213125006
$ws.Range("L29").Value = 22894.08909150676
$ws.Range("M29").Value = 803.2754764682949

$ws.Range("C30").Value = 1420.4325771695444
$ws.Range("E30").Value = 408.6434465609806
$ws.Range("F30").Value = 92.44933318643962
$ws.Range("G30").Value = 737.2840620933006
$ws.Range("K30").Value = 0.027773201774999996
$ws.Range("L30").Value = 20184.601146307563
$ws.Range("M30").Value = 560.5910003842962

$ws.Range("C31").Value = 1399.6031274815796
$ws.Range("E31").Value = 370.0282563280702
$ws.Range("F31").Value = 56.27396017370953
$ws.Range("G31").Value = 837.786013284828
$ws.Range("K31").Value = 0.0342295942225
$ws.Range("L31").Value = 20431.203681213374
$ws.Range("M31").Value = 699.351811485182

$r = $ws.Range("A26:N31")
$r.Select()